$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 7005.25
$ws.Cells.Item(20, 9).Value = 6007
$ws.Cells.Item(20, 10).Value = 10000
$ws.Cells.Item(20, 11).Value = 6007
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 13).Value = -5777
$ws.Cells.Item(20, 14).Value = -10460

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(35, 8).Value = 7005.25
$ws.Cells.Item(35, 9).Value = 6007
$ws.Cells.Item(35, 10).Value = 10000
$ws.Cells.Item(35, 11).Value = 6007
$ws.Cells.Item(35, 12).Value = 10000
$ws.Cells.Item(35, 13).Value = -5628
$ws.Cells.Item(35, 14).Value = -10758

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2867.8462
$ws.Cells.Item(62, 9).Value = 2426.889
$ws.Cells.Item(62, 10).Value = 3860
$ws.Cells.Item(62, 11).Value = 2426.889
$ws.Cells.Item(62, 12).Value = 3860
$ws.Cells.Item(62, 13).Value = -1802.889
$ws.Cells.Item(62, 14).Value = -5108

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2867.8462
$ws.Cells.Item(65, 9).Value = 2426.889
$ws.Cells.Item(65, 10).Value = 3860
$ws.Cells.Item(65, 11).Value = 12134.445
$ws.Cells.Item(65, 12).Value = 19300
$ws.Cells.Item(65, 13).Value = -9014.445
$ws.Cells.Item(65, 14).Value = -25540

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18523.512
$ws.Cells.Item(32, 10).Value = 23749.125
$ws.Cells.Item(32, 12).Value = 23749.125
$ws.Cells.Item(32, 14).Value = -24323.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3189.4856
$ws.Cells.Item(31, 9).Value = 3981
$ws.Cells.Item(31, 10).Value = 2826.7083
$ws.Cells.Item(31, 11).Value = 3981
$ws.Cells.Item(31, 12).Value = 2826.7083
$ws.Cells.Item(31, 13).Value = -3686
$ws.Cells.Item(31, 14).Value = -3416.7083

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3189.4856
$ws.Cells.Item(34, 9).Value = 3981
$ws.Cells.Item(34, 10).Value = 2826.7083
$ws.Cells.Item(34, 11).Value = 3981
$ws.Cells.Item(34, 12).Value = 2826.7083
$ws.Cells.Item(34, 13).Value = -3779
$ws.Cells.Item(34, 14).Value = -3230.7083

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 47621004
$ws.Cells.Item(58, 9).Value = 52633532
$ws.Cells.Item(58, 10).Value = 2000
$ws.Cells.Item(58, 11).Value = 52633532
$ws.Cells.Item(58, 12).Value = 2000
$ws.Cells.Item(58, 13).Value = -52633329
$ws.Cells.Item(58, 14).Value = -2406

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3601.2
$ws.Cells.Item(62, 9).Value = 3333.3333
$ws.Cells.Item(62, 11).Value = 3333.3333
$ws.Cells.Item(62, 13).Value = -2709.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 3601.2
$ws.Cells.Item(65, 9).Value = 3333.3333
$ws.Cells.Item(65, 11).Value = 16666.6665
$ws.Cells.Item(65, 13).Value = -13546.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2606.75
$ws.Cells.Item(86, 9).Value = 1877.7778
$ws.Cells.Item(86, 10).Value = 3544
$ws.Cells.Item(86, 11).Value = 1877.7778
$ws.Cells.Item(86, 12).Value = 3544
$ws.Cells.Item(86, 13).Value = -754.7778000000001
$ws.Cells.Item(86, 14).Value = -5790

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 2606.75
$ws.Cells.Item(89, 9).Value = 1877.7778
$ws.Cells.Item(89, 10).Value = 3544
$ws.Cells.Item(89, 11).Value = 9388.889000000001
$ws.Cells.Item(89, 12).Value = 17720
$ws.Cells.Item(89, 13).Value = -3772.889000000001
$ws.Cells.Item(89, 14).Value = -28952

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 537.2222
$ws.Cells.Item(107, 9).Value = 537.2222
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 537.2222
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 1382.7778
$ws.Cells.Item(107, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2472.077
$ws.Cells.Item(122, 9).Value = 2264.6
$ws.Cells.Item(122, 10).Value = 2601.75
$ws.Cells.Item(122, 11).Value = 6793.799999999999
$ws.Cells.Item(122, 12).Value = 7805.25
$ws.Cells.Item(122, 13).Value = -4343.799999999999
$ws.Cells.Item(122, 14).Value = -12705.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 36379.758
$ws.Cells.Item(132, 9).Value = 1760.04
$ws.Cells.Item(132, 10).Value = 252753
$ws.Cells.Item(132, 11).Value = 5280.12
$ws.Cells.Item(132, 12).Value = 758259
$ws.Cells.Item(132, 13).Value = -2750.12
$ws.Cells.Item(132, 14).Value = -763319

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 48423.957
$ws.Cells.Item(134, 9).Value = 1394.3529
$ws.Cells.Item(134, 10).Value = 181674.5
$ws.Cells.Item(134, 11).Value = 4183.0587
$ws.Cells.Item(134, 12).Value = 545023.5
$ws.Cells.Item(134, 13).Value = -1648.0587
$ws.Cells.Item(134, 14).Value = -550093.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 47621004
$ws.Cells.Item(136, 9).Value = 52633532
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 157900596
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = -157898046
$ws.Cells.Item(136, 14).Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 333300
$ws.Cells.Item(56, 9).Value = 333300
$ws.Cells.Item(56, 11).Value = 333300
$ws.Cells.Item(56, 13).Value = -332770

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 985.2093
$ws.Cells.Item(131, 10).Value = 1012.6
$ws.Cells.Item(131, 12).Value = 3037.8
$ws.Cells.Item(131, 14).Value = -13117.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2430.8965
$ws.Cells.Item(132, 9).Value = 2000.7
$ws.Cells.Item(132, 10).Value = 2657.3157
$ws.Cells.Item(132, 11).Value = 18006.3
$ws.Cells.Item(132, 12).Value = 23915.8413
$ws.Cells.Item(132, 13).Value = -15476.3
$ws.Cells.Item(132, 14).Value = -28975.8413

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 12083.071
$ws.Cells.Item(141, 9).Value = 4186
$ws.Cells.Item(141, 10).Value = 16470.334
$ws.Cells.Item(141, 11).Value = 12558
$ws.Cells.Item(141, 12).Value = 49411.00199999999
$ws.Cells.Item(141, 13).Value = -7378
$ws.Cells.Item(141, 14).Value = -59771.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 21414.285
$ws.Cells.Item(29, 9).Value = 4950
$ws.Cells.Item(29, 10).Value = 28000
$ws.Cells.Item(29, 11).Value = 4950
$ws.Cells.Item(29, 12).Value = 28000
$ws.Cells.Item(29, 13).Value = -4660
$ws.Cells.Item(29, 14).Value = -28580

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1557.5625
$ws.Cells.Item(113, 9).Value = 1378.9
$ws.Cells.Item(113, 10).Value = 1855.3334
$ws.Cells.Item(113, 11).Value = 1378.9
$ws.Cells.Item(113, 12).Value = 1855.3334
$ws.Cells.Item(113, 13).Value = 791.0999999999999
$ws.Cells.Item(113, 14).Value = -6195.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10872316
$ws.Cells.Item(7, 9).Value = 16668585
$ws.Cells.Item(7, 10).Value = 4311.375
$ws.Cells.Item(7, 11).Value = 16668585
$ws.Cells.Item(7, 12).Value = 4311.375
$ws.Cells.Item(7, 13).Value = -16668473
$ws.Cells.Item(7, 14).Value = -4535.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2287.0476
$ws.Cells.Item(61, 9).Value = 2495.1875
$ws.Cells.Item(61, 10).Value = 1621
$ws.Cells.Item(61, 11).Value = 2495.1875
$ws.Cells.Item(61, 12).Value = 1621
$ws.Cells.Item(61, 13).Value = -2293.1875
$ws.Cells.Item(61, 14).Value = -2025

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1481.25
$ws.Cells.Item(68, 10).Value = 1580
$ws.Cells.Item(68, 12).Value = 1580
$ws.Cells.Item(68, 14).Value = -3078

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 1481.25
$ws.Cells.Item(71, 10).Value = 1580
$ws.Cells.Item(71, 12).Value = 7900
$ws.Cells.Item(71, 14).Value = -15388

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1343.5652
$ws.Cells.Item(100, 9).Value = 1111.8823
$ws.Cells.Item(100, 11).Value = 1111.8823
$ws.Cells.Item(100, 13).Value = -570.8823

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 2287.0476
$ws.Cells.Item(113, 9).Value = 2495.1875
$ws.Cells.Item(113, 10).Value = 1621
$ws.Cells.Item(113, 11).Value = 2495.1875
$ws.Cells.Item(113, 12).Value = 1621
$ws.Cells.Item(113, 13).Value = -325.1875
$ws.Cells.Item(113, 14).Value = -5961

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3235.2456
$ws.Cells.Item(122, 9).Value = 2758.2
$ws.Cells.Item(122, 10).Value = 3493.1082
$ws.Cells.Item(122, 11).Value = 8274.599999999999
$ws.Cells.Item(122, 12).Value = 10479.3246
$ws.Cells.Item(122, 13).Value = -5824.599999999999
$ws.Cells.Item(122, 14).Value = -15379.3246

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 10872316
$ws.Cells.Item(126, 9).Value = 16668585
$ws.Cells.Item(126, 10).Value = 4311.375
$ws.Cells.Item(126, 11).Value = 50005755
$ws.Cells.Item(126, 12).Value = 12934.125
$ws.Cells.Item(126, 13).Value = -50003285
$ws.Cells.Item(126, 14).Value = -17874.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1974.963
$ws.Cells.Item(81, 9).Value = 1551.2142
$ws.Cells.Item(81, 10).Value = 2431.3076
$ws.Cells.Item(81, 11).Value = 3102.4284
$ws.Cells.Item(81, 12).Value = 4862.6152
$ws.Cells.Item(81, 13).Value = -2041.4284
$ws.Cells.Item(81, 14).Value = -6984.6152

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 1974.963
$ws.Cells.Item(84, 9).Value = 1551.2142
$ws.Cells.Item(84, 10).Value = 2431.3076
$ws.Cells.Item(84, 11).Value = 15512.142
$ws.Cells.Item(84, 12).Value = 24313.076
$ws.Cells.Item(84, 13).Value = -10208.142
$ws.Cells.Item(84, 14).Value = -34921.076
